# Chokhatauri disability_prevalence.xlsx update
# Switches the sheet from "Internally Displaced Disability Persons receiving
# Social Package" data (with confidential "..." placeholders) to the new
# "Unified database of targeted social assistance program" data, with two
# fully populated data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 1: title - new wording, now merged across A1:I1, centered + wrapped
# ---------------------------------------------------------------------
if ($ws.Range("A5").MergeCells) { $ws.Range("A5").UnMerge() }

$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Chokhatauri Municipality"
$ws.Range("A1:I1").Merge()
$r1 = $ws.Range("A1:I1")
$r1.HorizontalAlignment = -4108   # xlCenter
$r1.VerticalAlignment = -4108     # xlCenter
$r1.WrapText = $true
$r1.Font.Name = "Arial"
$r1.Font.Size = 11
$r1.Font.Bold = $true
$r1.Borders.Item(7).LineStyle = -4142
$r1.Borders.Item(8).LineStyle = -4142
$r1.Borders.Item(9).LineStyle = -4142
$r1.Borders.Item(10).LineStyle = -4142
$r1.Interior.Pattern = -4142      # xlNone - no fill
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------------
# Row 2: unchanged text/formatting - "(End of year, persons)"
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "(End of year, persons)"
$ws.Rows.Item(2).RowHeight = 14.5

# ---------------------------------------------------------------------
# Row 3: year headers stay the same; A3 (blank) switches to Sylfaen font
# ---------------------------------------------------------------------
$a3 = $ws.Range("A3")
$a3.Font.Name = "Sylfaen"
$a3.Font.Size = 11

# ---------------------------------------------------------------------
# Row 4: "family with disabilities Persons" data row
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$data4 = @(692, 659, 620, 631, 608, 1472, 592, 569)
$cols = @("B","C","D","E","F","G","H","I")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value = $data4[$i]
}

$row4 = $ws.Range("A4:I4")
$row4.Font.Name = "Arial"
$row4.Font.Size = 10
$row4.Font.Bold = $false
$row4.Font.Underline = -4142
$row4.Interior.Pattern = 1
$row4.Interior.ThemeColor = 2
$row4.Interior.TintAndShade = 0
$ws.Range("B4:I4").NumberFormat = "# ##0"
$ws.Range("B4:I4").HorizontalAlignment = 1   # xlGeneral
$ws.Range("A4").HorizontalAlignment = -4131  # xlLeft
$ws.Range("A4").VerticalAlignment = -4108    # xlCenter
$ws.Range("A4").WrapText = $true
$row4.Borders.Item(9).LineStyle = -4142      # clear bottom border
$row4.Borders.Item(8).LineStyle = 1          # keep/ensure thin top border
$row4.Borders.Item(8).Weight = 2
$ws.Rows.Item(4).RowHeight = 24.75

# ---------------------------------------------------------------------
# Row 5: "disabilities Persons" data row (was the merged "Note" row)
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$data5 = @(777, 743, 703, 716, 696, 1633, 677, 646)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value = $data5[$i]
}

$row5 = $ws.Range("A5:I5")
$row5.Font.Name = "Arial"
$row5.Font.Size = 10
$row5.Font.Bold = $false
$row5.Font.Underline = -4142
$row5.Interior.Pattern = 1
$row5.Interior.ThemeColor = 2
$row5.Interior.TintAndShade = 0
$ws.Range("B5:I5").NumberFormat = "# ##0"
$ws.Range("B5:I5").HorizontalAlignment = 1   # xlGeneral
$ws.Range("A5").HorizontalAlignment = -4131  # xlLeft
$ws.Range("A5").VerticalAlignment = -4108    # xlCenter
$ws.Range("A5").WrapText = $true
$row5.Borders.Item(7).LineStyle = -4142
$row5.Borders.Item(8).LineStyle = -4142
$row5.Borders.Item(10).LineStyle = -4142
$row5.Borders.Item(9).LineStyle = 1          # thin bottom border (closes table)
$row5.Borders.Item(9).Weight = 2
$ws.Rows.Item(5).RowHeight = 21

# ---------------------------------------------------------------------
# Row 6: "Source: ..." note, now merged A6:H6, plain (not bold/underlined)
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "Source: Ministry of Internally Displaced Persons from the Occupied Territories, Labour, Health and Social Affairs of Georgia."
$ws.Range("A6:H6").Merge()
$row6 = $ws.Range("A6:H6")
$row6.Font.Name = "Arial"
$row6.Font.Size = 9
$row6.Font.Bold = $false
$row6.Font.Underline = -4142
$row6.Interior.Pattern = 1
$row6.Interior.ThemeColor = 2
$row6.Interior.TintAndShade = 0
$row6.HorizontalAlignment = -4131   # xlLeft
$row6.VerticalAlignment = -4108     # xlCenter
$row6.WrapText = $true
$row6.Borders.Item(7).LineStyle = -4142
$row6.Borders.Item(9).LineStyle = -4142
$row6.Borders.Item(10).LineStyle = -4142
$ws.Range("B6:H6").Borders.Item(8).LineStyle = 1
$ws.Range("B6:H6").Borders.Item(8).Weight = 2
$ws.Range("A6").Borders.Item(8).LineStyle = -4142
$ws.Rows.Item(6).RowHeight = 27.75

# ---------------------------------------------------------------------
# Column widths: only column A keeps a custom width now
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.8164
for ($c = 2; $c -le 16; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 8.43
}

Write-Output "Chokhatauri disability_prevalence sheet updated"
